$wb = $excel.ActiveWorkbook

# Overview sheet: update status for fa496268 file (row 3) from
# "Ready for handoff" to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status + handback datetime for fa496268 file (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-01-28 05:29:49"

# de-de sheet: update status + handback datetime for fa496268 file (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-01-28 05:30:10"
